# trash.pantallas.checklist.xlsx - mark additional screens as "maquetada"/"implementada"
# (commit: "cambios db, marco avances de programacion, y subo nuevo DER")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Mark progress ("x") on rows whose column B/C checkboxes were not yet set.
$ws.Range("B3").Value = "x"
$ws.Range("B7").Value = "x"
$ws.Range("B17").Value = "x"
$ws.Range("B23").Value = "x"
$ws.Range("C23").Value = "x"
$ws.Range("B24").Value = "x"
$ws.Range("C24").Value = "x"
$ws.Range("B25").Value = "x"
$ws.Range("C25").Value = "x"

# Update the view state: scroll the window so row 13 is at the top, and
# move the active selection to D6.
[void]$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
[void]$ws.Range("D6").Select()
